$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "[강화학습] Reinforcement learning for Process Control (Python, torch)"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/%EA%B0%95%ED%99%94%ED%95%99%EC%8A%B5-%EA%B0%95%ED%99%94%ED%95%99%EC%8A%B5%EC%9C%BC%EB%A1%9C-%EA%B3%B5%EC%A0%95-%EC%9E%90%EB%8F%99-%EC%A0%9C%EC%96%B4-%EC%8B%9C%EB%AE%AC%EB%A0%88%EC%9D%B4%EC%85%98-%ED%95%B4%EB%B3%B4%EA%B8%B0"

$ws.Range("D32").Value = "TensorFlow 모델 서빙 방법 (feat. ChatGPT)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/411"
